# Updated cryptos list on Tue Jan 30 17:33:47 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must remain stored as TEXT even when it looks
# like a number (e.g. "310.23"), without leaving a custom number-format
# applied to the cell afterwards (mirrors the original inlineStr cells,
# which carry no explicit style).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "43.413.79"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.378.45"
$ws.Range("E3").Value = "  +3.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "310.23"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "104.59"
$ws.Range("E6").Value = "  +3.21%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.526"
$ws.Range("E7").Value = "  -1.66%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.55%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "36.28"

# Row 11 - OKB
Set-TextValue $ws.Range("D11") "52.76"
$ws.Range("E11").Value = "  +1.05%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -0.73%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.78%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "7.00"
$ws.Range("E14").Value = "  +0.25%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "2.747.33"
$ws.Range("E15").Value = "  +3.35%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "15.76"
$ws.Range("E16").Value = "  +6.04%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.379.92"
$ws.Range("E17").Value = "  +3.92%  "

# Row 18 - Polygon
Set-TextValue $ws.Range("D18") "0.820"
$ws.Range("E18").Value = "  +1.87%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "43.395.28"
$ws.Range("E19").Value = "  +0.78%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D20") "12.00"
$ws.Range("E20").Value = "  -4.24%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.05%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +3.68%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "68.53"
$ws.Range("E23").Value = "  +0.49%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "242.33"
$ws.Range("E24").Value = "  +0.73%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  +1.72%  "

# Row 26 - PancakeSwap
Set-TextValue $ws.Range("D26") "2.62"
$ws.Range("E26").Value = "  -0.38%  "

# Row 27 - Dai
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.39%  "

# Row 28 - EthereumClassic
Set-TextValue $ws.Range("D28") "26.18"
$ws.Range("E28").Value = "  +8.53%  "

# Row 29 - Toncoin
Set-TextValue $ws.Range("D29") "2.27"
$ws.Range("E29").Value = "  +5.04%  "

# Row 30 - InjectiveProtocol
Set-TextValue $ws.Range("D30") "36.66"
$ws.Range("E30").Value = "  -4.41%  "

# Row 31 - Cosmos
$ws.Range("E31").Value = "  -0.49%  "

# Row 32 - Monero
Set-TextValue $ws.Range("D32") "161.53"
$ws.Range("E32").Value = "  -2.47%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -0.99%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  +0.08%  "

# Row 35 - Celestia
Set-TextValue $ws.Range("D35") "18.32"
$ws.Range("E35").Value = "  +2.60%  "

# Row 36 - LidoDAOToken
Set-TextValue $ws.Range("D36") "3.15"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37 - WEMIXToken
$ws.Range("E37").Value = "  +6.54%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  +0.25%  "

# Row 39 - RenderToken
Set-TextValue $ws.Range("D39") "4.68"
$ws.Range("E39").Value = "  +11.41%  "

# Row 40 - ARBITRUM
Set-TextValue $ws.Range("D40") "1.96"
$ws.Range("E40").Value = "  +5.92%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +0.57%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  -1.16%  "

# Row 43 - ApeXProtocol
Set-TextValue $ws.Range("D43") "2.44"
$ws.Range("E43").Value = "  +5.43%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "20.50"
$ws.Range("E44").Value = "  +3.55%  "

# Row 45 - Maker
Set-TextValue $ws.Range("D45") "2.011.41"
$ws.Range("E45").Value = "  +2.22%  "

# Row 46 / 47 - NEARProtocol and VeChain swap list positions
# (NEARProtocol moves into row 46, VeChain moves into row 47)
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D46") "3.21"
$ws.Range("E46").Value = "  +5.77%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0291"
$ws.Range("E47").Value = "  +0.60%  "

# Row 48 - FraxShare
Set-TextValue $ws.Range("D48") "10.51"
$ws.Range("E48").Value = "  +6.89%  "

# Row 49 - MultiversX
Set-TextValue $ws.Range("D49") "58.03"
$ws.Range("E49").Value = "  +5.22%  "

# Row 50 - HuobiToken
Set-TextValue $ws.Range("D50") "2.92"
$ws.Range("E50").Value = "  -3.30%  "

# Row 51 - RocketPoolETH
Set-TextValue $ws.Range("D51") "2.581.60"
$ws.Range("E51").Value = "  +2.15%  "
